# "Înput Data Base Case"
# Insert a new leading index column (A) in front of the existing
# Technology / Invested capacity (MW) table, shifting the old A/B
# columns to B/C, and fill in the base-case input values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing Technology (A) / Invested capacity (B) columns one
# position to the right, creating a new, currently-empty column A.
$ws.Range("A1").EntireColumn.Insert()

# Populate the new column A with a 0-based row index for each data row
# (rows 2-17 -> 0-15). The header row (row 1) gets no value in column A.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Give the new index cells the same bold/centered/bordered formatting
# already used by the header cells (B1/C1).
$ws.Range("B1").Copy()
$ws.Range("A2:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the "base case" invested capacity (MW) values in column C.
$ws.Range("C7").Value = 400
$ws.Range("C8").Value = 200
$ws.Range("C9").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 400
$ws.Range("C17").Value = 0
